$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Move the "_GoBack" bookmark from its old location (after the "Sage MAS "
#    run near "90 programs") to right after "Charles A. Rolke" in the title
#    paragraph. The host models "_GoBack" as a singleton bookmark, so simply
#    adding a new one removes the old one automatically.
#
#    A collapsed bookmark placed exactly at the last text position of a
#    paragraph (i.e. right before its paragraph mark) gets normalised to the
#    paragraph start on save, so we insert a temporary 2-character placeholder
#    after "Charles A. Rolke", anchor the bookmark just before it (now a
#    mid-paragraph position), and then delete the placeholder again. The
#    collapsed bookmark stays put, ending up correctly between the run and
#    the paragraph mark.
# ---------------------------------------------------------------------------
$rTitle = $d.Content
$okTitle = $rTitle.Find.Execute("Charles A. Rolke", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okTitle) { throw "Could not find 'Charles A. Rolke'" }
$rTitle.Collapse(0)
$anchorStart = $rTitle.Start

$rPlaceholder = $d.Range($anchorStart, $anchorStart)
$rPlaceholder.InsertAfter("ZZ")

$rBookmark = $d.Range($anchorStart, $anchorStart)
$d.Bookmarks.Add("_GoBack", $rBookmark)

$rPlaceholder2 = $d.Range($anchorStart, $anchorStart + 2)
$rPlaceholder2.Delete()

# ---------------------------------------------------------------------------
# 2. Tighten the bottom paragraph-border spacing (w:space 3 -> 0) on the two
#    contact-info paragraphs (phone/email line and the blank line under it).
# ---------------------------------------------------------------------------
$pPhone = $d.Paragraphs(3)
$pPhone.Borders.DistanceFromBottom = 0

$pBlank = $d.Paragraphs(4)
$pBlank.Borders.DistanceFromBottom = 0

# ---------------------------------------------------------------------------
# 3. Split "(508) 250-3076 / crolke2000@yahoo.com" into a plain run, a
#    hyperlinked e-mail address, and a new line with a web reference.
# ---------------------------------------------------------------------------
$rReplace = $d.Content
$okReplace = $rReplace.Find.Execute( `
    "(508) 250-3076 / crolke2000@yahoo.com", $false, $false, $false, $false, $false, $true, 1, $false, `
    "(508) 250-3076 / crolke2000@yahoo.com^lonline at https://github.com/carolke/Resume", 2)
if (-not $okReplace) { throw "Could not find/replace the phone/e-mail line" }

$rEmail = $d.Content
$okEmail = $rEmail.Find.Execute("crolke2000@yahoo.com", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okEmail) { throw "Could not re-find the e-mail address to hyperlink" }
$d.Hyperlinks.Add($rEmail, "mailto:crolke2000@yahoo.com", "", "", "crolke2000@yahoo.com") | Out-Null

# ---------------------------------------------------------------------------
# 4. Drop the extra empty paragraph that trailed the "Associate of Science"
#    line right before the section break.
# ---------------------------------------------------------------------------
$countParas = $d.Paragraphs.Count
$pLast = $d.Paragraphs($countParas)
$pPrev = $d.Paragraphs($countParas - 1)
$rTrim = $d.Range($pPrev.Range.End - 1, $pLast.Range.End)
$rTrim.Delete()

Write-Output "done"
